$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.730.19"
$ws.Range("E2").Value = "  +1.94%  "

$ws.Range("D3").Value = "1.638.12"
$ws.Range("E3").Value = "  +2.05%  "

$ws.Range("E4").Value = "  -0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.79"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.14%  "

$ws.Range("E6").Value = "  +1.92%  "

$ws.Range("E7").Value = "  -0.11%  "

$ws.Range("E8").Value = "  +1.27%  "

$ws.Range("E9").Value = "  +1.98%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.09"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.94%  "

$ws.Range("E11").Value = "  +2.80%  "

$ws.Range("D12").Value = "1.866.76"
$ws.Range("E12").Value = "  +2.07%  "

$ws.Range("D13").Value = "1.631.72"
$ws.Range("E13").Value = "  +1.60%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.08"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.70%  "

$ws.Range("E15").Value = "  +2.88%  "

$ws.Range("D16").Value = "26.742.68"
$ws.Range("E16").Value = "  +2.10%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.12"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.32%  "

$ws.Range("D18").Value = "0.0$([char]0x2083)" + "0742"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "209.05"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.24%  "

$ws.Range("E20").Value = "  -0.22%  "

$ws.Range("E21").Value = "  +1.15%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.41"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.55%  "

$ws.Range("E23").Value = "  +2.54%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.93"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.82%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "146.51"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.73%  "

$ws.Range("E26").Value = "  -0.18%  "

$ws.Range("E27").Value = "  -0.31%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.78"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.72%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.41"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.56%  "

$ws.Range("E30").Value = "  +5.81%  "

$ws.Range("E31").Value = "  -0.33%  "

$ws.Range("E32").Value = "  +1.16%  "

$ws.Range("E33").Value = "  +1.59%  "

$ws.Range("E34").Value = "  +1.03%  "

$ws.Range("E35").Value = "  +0.52%  "

$ws.Range("D36").Value = "1.169.27"
$ws.Range("E36").Value = "  +0.43%  "

$ws.Range("E37").Value = "  -0.13%  "

$ws.Range("E38").Value = "  +3.31%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.505"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.90%  "

$ws.Range("E41").Value = "  +0.35%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.796"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.43%  "

$ws.Range("E43").Value = "  +1.62%  "

$ws.Range("D44").Value = "1.776.36"
$ws.Range("E44").Value = "  +2.06%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "92.52"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.92%  "

$ws.Range("E46").Value = "  +2.78%  "

$ws.Range("D47").Value = "0.0$([char]0x2086)" + "0104"
$ws.Range("E47").Value = "  +7.12%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "54.76"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.35%  "

$ws.Range("E49").Value = "  +1.47%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.409"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.58%  "

$ws.Range("E51").Value = "  +4.22%  "
